$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Overview')
$ws.Range('A2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.md'
$ws.Range('D2').Value = '2016-06-17 06:06:08'
$ws.Range('A3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.md'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('D3').Value = '2016-04-17 06:04:38'
$ws = $wb.Worksheets.Item('zh-cn')
$ws.Range('A2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.md'
$ws.Range('D2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.zh-cn.xlf'
$ws.Range('E2').Value = '2016-03-17 06:06:00'
$ws.Range('F2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.md'
$ws.Range('G2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.zh-cn.xlf'
$ws.Range('H2').Value = '2016-03-17 06:06:40'
$ws.Range('A3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.md'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('D3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.9695ea9775b82014281b12a7a9e1234ace738f2e.zh-cn.xlf'
$ws.Range('E3').Value = '2016-03-17 06:04:30'
$ws.Range('F3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.md'
$ws.Range('G3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.9695ea9775b82014281b12a7a9e1234ace738f2e.zh-cn.xlf'
$ws = $wb.Worksheets.Item('de-de')
$ws.Range('A2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.md'
$ws.Range('D2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.de-de.xlf'
$ws.Range('E2').Value = '2016-03-17 06:06:08'
$ws.Range('F2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.md'
$ws.Range('G2').Value = 'ae18194b-63eb-4d66-9d9a-8155b62d5b35.3419fc0ff14cca1722c8dad84163dd260b0ab66d.de-de.xlf'
$ws.Range('H2').Value = '2016-03-17 06:06:53'
$ws.Range('A3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.md'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('D3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.9695ea9775b82014281b12a7a9e1234ace738f2e.de-de.xlf'
$ws.Range('E3').Value = '2016-03-17 06:04:38'
$ws.Range('F3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.md'
$ws.Range('G3').Value = 'f0d6173c-d2fa-4d72-81a9-30eed75f0cff.9695ea9775b82014281b12a7a9e1234ace738f2e.de-de.xlf'
